$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data values (Sun Dec 18 19:09:33 UTC 2022 symbol list refresh).
# Numeric-looking values must keep their exact text representation (e.g. trailing zeros,
# 4-decimal formatting), so those cells are forced to Text format before assignment.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.85"
$ws.Range("E2").Value = "1BNBBNB"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "19"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.73"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "19"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.521"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "19"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05627"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "19"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.469"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "19"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8069"
$ws.Range("E7").Value = "6MXTokenMX"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "19"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.034"
$ws.Range("E8").Value = "7FTXTokenFTT"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "19"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1430"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "19"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07320"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "19"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03104"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "19"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02913"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "19"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09267"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "19"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001663"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "19"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.216"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "19"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04739"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "19"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.01167"
$ws.Range("E17").Value = "16OneONEBestin24h"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "19"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006427"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "19"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005070"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "19"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001053"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "19"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001501"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "19"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.981"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "19"
$ws.Range("B23").Value = "GateToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.380"
$ws.Range("E23").Value = "22GateTokenGT"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "19"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.105"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "19"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3268"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "19"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1257"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "19"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003302"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "19"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "19"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "19"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "19"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "19"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "19"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "19"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "19"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "19"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "19"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "19"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "19"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "19"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04146"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "19"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007084"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "19"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1041"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "19"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003302"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "19"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008711"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "19"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005651"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "19"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "19"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6804"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "19"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.01633"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "19"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "19"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "19"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "19"
